# Sync attendance_reports: rotate the leading "System" entry in the
# "Recorded By" column (G) to the end of the comma-separated list,
# e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".
# Cells that don't begin with "System, " (e.g. a bare "System", or an
# entry that doesn't start with System at all) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$prefix = "System, "

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -ne $val -and $val.StartsWith($prefix)) {
        $rest = $val.Substring($prefix.Length)
        $cell.Value2 = $rest + ", System"
    }
}
